# "Atualizacao de status do leito"
#
# The "grants por usuario" sheet repeats one GRANT statement (column A,
# rows 100-151) that is concatenated with a username (column B) and a
# trailing ";" (column C) via the existing shared formula in column D
# (=A&" "&B&" "&C). Swap the GRANT text for the new SEQUENCE grant; the
# formula results in column D recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("grants por usuario")

$ws.Range("A100:A151").Value = "GRANT ALL ON SEQUENCE integracao.sq_hstr_ocpa_leito_status TO"
